$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36, shifting existing rows 36-39 down to 37-40.
$ws.Rows.Item(36).Insert()

# Fill in the new row 36 with the new weekly data record.
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 45127
$ws.Cells.Item(36, 4).NumberFormat = $ws.Cells.Item(37, 4).NumberFormat
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 100112013
$ws.Cells.Item(36, 7).Value = "Alcachofa"
$ws.Cells.Item(36, 8).Value = "Madrigal"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 70
$ws.Cells.Item(36, 11).Value = 19000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 19357
$ws.Cells.Item(36, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(36, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(36, 16).Value = 484
$ws.Cells.Item(36, 17).Value = 40
$ws.Cells.Item(36, 18).Value = "Hortaliza"

# Row 38 (previously row 37) had its Origen value changed.
$ws.Cells.Item(38, 15).Value = "Región de Coquimbo"
